$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1").EntireColumn.Delete()
foreach ($hl in $ws.Hyperlinks) {
    Write-Host "after:" $hl.Range.Address()
}
